$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "20250526_101139"
$ws.Range("B3").Value = "2025-05-26 10:11:39"
$ws.Range("C3").Value = "Isabelle Roux"
$ws.Range("D3").Value = "{'chantier': 'Métallerie - Table Métal 08', 'urgence': 'Normal', 'date_souhaitee': '2025-05-26', 'produits': {2938475610: {'produit': 'Marteau 500g', 'quantite': 4, 'emplacement': 'Atelier B'}, 6574839202: {'produit': 'Vis 8x60mm', 'quantite': 1, 'emplacement': 'Stockage'}, 1928374650: {'produit': 'Clé à molette', 'quantite': 4, 'emplacement': 'Atelier B'}}}"
$ws.Range("E3").Value = "pvc04`n"
$ws.Range("F3").Value = "Refusée"
$ws.Range("G3").Value = "2025-05-26 12:10:32"
$ws.Range("H3").Value = "Magasinier"
$ws.Range("I3").Value = "a"

# Row 4
$ws.Range("A4").Value = "20250526_112800"
$ws.Range("B4").Value = "2025-05-26 11:28:00"
$ws.Range("C4").Value = "Sophie Leroy"
$ws.Range("D4").Value = "{'chantier': 'PVC - Table PVC 04', 'urgence': 'Normal', 'date_souhaitee': '2025-05-26', 'produits': {1928374650: {'produit': 'Clé à molette', 'quantite': 5, 'emplacement': 'Atelier B'}}}"
$ws.Range("E4").Value = "r"
$ws.Range("F4").Value = "En attente"

# Row 5
$ws.Range("A5").Value = "20250526_113412"
$ws.Range("B5").Value = "2025-05-26 11:34:12"
$ws.Range("C5").Value = "Marie Martin"
$ws.Range("D5").Value = "{'chantier': 'Aluminium - Table Aluminium 02', 'urgence': 'Normal', 'date_souhaitee': '2025-05-26', 'produits': {2938475610: {'produit': 'Marteau 500g', 'quantite': 4, 'emplacement': 'Atelier B'}, 8473926150: {'produit': 'Tournevis cruciforme', 'quantite': 4, 'emplacement': 'Atelier A'}, 9182736450: {'produit': 'Perceuse sans fil', 'quantite': 2, 'emplacement': 'Atelier A'}}}"
$ws.Range("E5").Value = "h"
$ws.Range("F5").Value = "En attente"

# Row 6
$ws.Range("A6").Value = "20250526_122347"
$ws.Range("B6").Value = "2025-05-26 12:23:47"
$ws.Range("C6").Value = "Marie Martin"
$ws.Range("D6").Value = "{'chantier': 'Aluminium - Table Aluminium 02', 'urgence': 'Normal', 'date_souhaitee': '2025-05-26', 'produits': {7465839201: {'produit': 'Pince coupante', 'quantite': 40, 'emplacement': 'Atelier B'}}}"
$ws.Range("F6").Value = "En attente"

Write-Output "rows 3-6 written"
